$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing table (Tabela1) which currently spans A1:G5
$lo = $ws.ListObjects.Item(1)

# Add a new calculated column "Shelf" to the table (extends the table to A1:H5)
$newCol = $lo.ListColumns.Add()
$ws.Cells.Item(1, 8).Value = "Shelf"

# Give the new header cell the same formatting as the other header cells.
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Fill in the calculated-column formula for each data row, matching the
# per-row literal formulas (non shared) and cached string results that a
# real table calculated-column would produce.
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=RIGHT(Tabela1[[#This Row],[Y-Coor]], 1)"
}

# Apply the same yellow fill formatting used across the rest of the table
# to the new data cells.
$ws.Range("H2:H5").Interior.Color = 65535

# Update the active selection to match the authored state.
$null = $ws.Range("K5").Select()
